$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.788.28'
$ws.Range("E2").Value = '  +3.16%  '
$ws.Range("D3").Value = '3.678.22'
$ws.Range("E3").Value = '  +8.95%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '589.33'
$ws.Range("E5").Value = '  +1.47%  '
$ws.Range("D6").Value = '180.87'
$ws.Range("E6").Value = '  +1.06%  '
$ws.Range("D7").Value = '3.668.06'
$ws.Range("E7").Value = '  +8.98%  '
$ws.Range("D8").Value = '0.623'
$ws.Range("E8").Value = '  +4.89%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = '0.204'
$ws.Range("E10").Value = '  +3.12%  '
$ws.Range("D11").Value = '0.614'
$ws.Range("E11").Value = '  +4.41%  '
$ws.Range("D12").Value = '50.01'
$ws.Range("E12").Value = '  +3.57%  '
$ws.Range("D13").Value = '0.0000288'
$ws.Range("E13").Value = '  +1.48%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '4.259.31'
$ws.Range("E14").Value = '  +8.55%  '
$ws.Range("B15").Value = 'BitcoinCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D15").Value = '682.72'
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("D16").Value = '9.03'
$ws.Range("E16").Value = '  +4.94%  '
$ws.Range("D17").Value = '71.920.04'
$ws.Range("E17").Value = '  +3.30%  '
$ws.Range("D18").Value = '3.668.26'
$ws.Range("E18").Value = '  +8.42%  '
$ws.Range("E19").Value = '  +2.00%  '
$ws.Range("D20").Value = '18.30'
$ws.Range("E20").Value = '  +3.52%  '
$ws.Range("D21").Value = '11.69'
$ws.Range("E21").Value = '  +4.03%  '
$ws.Range("D22").Value = '0.945'
$ws.Range("E22").Value = '  +3.95%  '
$ws.Range("D23").Value = '6.22'
$ws.Range("E23").Value = '  +16.18%  '
$ws.Range("D24").Value = '17.86'
$ws.Range("E24").Value = '  +3.52%  '
$ws.Range("D25").Value = '103.72'
$ws.Range("E25").Value = '  +2.17%  '
$ws.Range("D26").Value = '4.03'
$ws.Range("E26").Value = '  +3.68%  '
$ws.Range("E27").Value = '  +5.47%  '
$ws.Range("D28").Value = '10.22'
$ws.Range("E28").Value = '  +5.57%  '
$ws.Range("D29").Value = '35.41'
$ws.Range("E29").Value = '  +5.62%  '
$ws.Range("D30").Value = '9.23'
$ws.Range("E30").Value = '  +5.78%  '
$ws.Range("D31").Value = '7.45'
$ws.Range("E31").Value = '  +8.04%  '
$ws.Range("D32").Value = '4.24'
$ws.Range("E32").Value = '  +10.52%  '
$ws.Range("D33").Value = '581.50'
$ws.Range("E33").Value = '  +4.81%  '
$ws.Range("D34").Value = '11.36'
$ws.Range("E34").Value = '  +2.91%  '
$ws.Range("E35").Value = '  +3.26%  '
$ws.Range("D36").Value = '60.25'
$ws.Range("E36").Value = '  +4.20%  '
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '3.741.84'
$ws.Range("E38").Value = '  +4.00%  '
$ws.Range("D39").Value = '0.144'
$ws.Range("E39").Value = '  +3.48%  '
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0776'
$ws.Range("E40").Value = '  +6.76%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = '35.69'
$ws.Range("E41").Value = '  +1.15%  '
$ws.Range("D42").Value = '3.47'
$ws.Range("E42").Value = '  +5.25%  '
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").Value = '2.80'
$ws.Range("E43").Value = '  +2.21%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0462'
$ws.Range("E44").Value = '  +8.62%  '
$ws.Range("D45").Value = '0.348'
$ws.Range("E45").Value = '  +3.62%  '
$ws.Range("D46").Value = '3.38'
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("D47").Value = '2.83'
$ws.Range("E47").Value = '  +6.20%  '
$ws.Range("E48").Value = '  +3.96%  '
$ws.Range("D49").Value = '1.45'
$ws.Range("E49").Value = '  +4.43%  '
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("D51").Value = '133.86'
$ws.Range("E51").Value = '  +2.69%  '
